$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match-detail columns (F:V) between two rows, leaving the
#     row-identity columns (A:E -- index/country/tournament/season/date)
#     untouched. Column range F..V is columns 6..22.
function Swap-RowDetails($r1, $r2) {
    $colStart = 6
    $colEnd = 22
    $vals1 = @()
    $vals2 = @()
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $vals1 += $ws.Cells.Item($r1, $c).Value2
        $vals2 += $ws.Cells.Item($r2, $c).Value2
    }
    $i = 0
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $ws.Cells.Item($r1, $c).Value2 = $vals2[$i]
        $ws.Cells.Item($r2, $c).Value2 = $vals1[$i]
        $i = $i + 1
    }
}

Swap-RowDetails 19 20
Swap-RowDetails 43 44
Swap-RowDetails 84 85
Swap-RowDetails 138 139

# --- Append two new match rows (147, 148) after the existing last row
#     (146), copying row 146's cell formatting first so the new rows
#     pick up the same styles (bordered/bold index in A, datetime format
#     in E) already used throughout the sheet.
$ws.Range("A146:V146").Copy()
$ws.Range("A147:V148").PasteSpecial(-4122)
$excel.CutCopyMode = $false

function Set-Row($r, $values) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range($cols[$i] + $r).Value2 = $values[$i]
    }
}

Set-Row 147 @(146, "serbia", "super-liga", "2023-2024", 45281.70833333334, `
    "Radnicki Nis", 1, "Zeleznicar Pancevo", 1, `
    1.71, "25/09/2023 05:12", 1.87, "21/12/2023 16:40", `
    3.52, "25/09/2023 05:12", 3.55, "21/12/2023 16:40", `
    4.19, "25/09/2023 05:12", 3.91, "21/12/2023 16:40", `
    "https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-zeleznicar-pancevo/QLI8f2FH/")

Set-Row 148 @(147, "serbia", "super-liga", "2023-2024", 45281.79166666666, `
    "Vozdovac", 3, "Vojvodina", 2, `
    3.21, "27/09/2023 06:42", 3.09, "21/12/2023 18:54", `
    3.2, "27/09/2023 06:42", 3.54, "21/12/2023 18:56", `
    2.04, "27/09/2023 06:42", 2.15, "21/12/2023 18:56", `
    "https://www.betexplorer.com/football/serbia/super-liga/fk-vozdovac-vojvodina/25ClaQ0n/")
